# Added trends for searches
# Adds two new rows (77 and 78) to Sheet1 describing new Redis keys used
# for tracking search trends (overall and by user), mirroring the existing
# "play" trend rows (71/72).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 77: overall search trends key
$ws.Range("A77").Value2 = "tinytape_<time code>_searches"
$ws.Range("B77").Value2 = "zset"
$ws.Range("C77").Value2 = "search count"
$ws.Range("D77").Value2 = "search term"
$ws.Range("E77").Value2 = "Search trends"

# Row 78: per-user search trends key
# (Note: column order below mirrors the shared-string insertion order of the
# original workbook edit.)
$ws.Range("C78").Value2 = "search cound"
$ws.Range("A78").Value2 = "tinytape_<time code>_searches_<username>"
$ws.Range("B78").Value2 = "zset"
$ws.Range("D78").Value2 = "search term"
$ws.Range("E78").Value2 = "Search trends by user"

# Update the active selection to reflect the new end of the used data,
# matching Excel's typical behavior of advancing the selection after the
# last filled row.
$ws.Range("E79").Select()
